# Regenerate orders with updated distance/sizes.
# Distance codes: D64 -> D69, D51 -> D55, D80 -> D86
# Size codes:     S30 -> S31  (S20 / S25 unchanged)
#
# These tokens show up embedded inside many different strings across the
# sheet (Condition, Filename_Left, Filename_Right, Distance, Size columns),
# so every string-valued cell in the used range is passed through the same
# substitution and written back if it changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -is [string]) {
            $newVal = $val
            # Use a placeholder so the D64->D69 pass can't collide with the
            # freshly-written D69 when later (it doesn't here, but keep the
            # substitution order collision-safe regardless).
            $newVal = $newVal.Replace("D64", "<<D69>>")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("<<D69>>", "D69")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
